# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the refreshed output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) — F column updates keyed by row number
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    2  = 568
    4  = 1234
    6  = 14010
    7  = 15303
    9  = 33
    10 = 45
    11 = 182
    18 = 70
    20 = 1182
    21 = 125
    23 = 5942
    24 = 954
    25 = 1077
    26 = 5503
    28 = 135
    29 = 84
    30 = 427
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# Sheet "全部类型" (sheet4) — F column updates keyed by row number
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    3  = 568
    5  = 1234
    7  = 14010
    8  = 15303
    10 = 33
    11 = 45
    12 = 182
    19 = 70
    21 = 1182
    22 = 125
    25 = 5942
    26 = 954
    27 = 1077
    28 = 5503
    30 = 135
    31 = 84
    32 = 427
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
